$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.945574641227722
$ws.Range("B1").Value = -1
$ws.Range("C1").Value = 2.230738162994385
$ws.Range("D1").Value = 1.422735214233398
$ws.Range("E1").Value = 1.153838992118835
